$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Selection moved to row 22 (whole row selected) ---
[void]$ws.Range("A22:XFD22").Select()

# --- Row 7: drop the now-unused O7/P7 style-only cells ---
$ws.Range("O7:P7").Clear()

# --- Row 9 (geo): new K9 example value, a single space placeholder ---
$ws.Range("K9").Value = " "

# --- Rows 12/13: clear the stray explicit 0 in column G (leave cell blank) ---
$ws.Range("G12").ClearContents()
$ws.Range("G13").ClearContents()

# --- Relabel the composite-family rows (18-22) to their new names/order.
#     Do these before row 14 so "COMPOSITE" is the last new shared string
#     introduced, matching the authored workbook's string table order. ---
$ws.Range("A18").Value = "composite-list"

$ws.Range("A19").Value = "map-of"

$ws.Range("A20").Value = "map-of-list"
$ws.Range("G20").ClearContents()

$ws.Range("A21").Value = "map-of-map"
$ws.Range("G21").ClearContents()

$ws.Range("A22").Value = "composite-map"
$ws.Range("E22").Value = 1
$ws.Range("F22").Value = 1
$ws.Range("H22").Value = 1

# --- Row 14: reserved -> COMPOSITE (new flag), clear stray 0 in column G ---
$ws.Range("A14").Value = "COMPOSITE"
$ws.Range("G14").ClearContents()
